$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 'maa://25390 (95.11), maa://24702 (94.79), maa://36681 (86.42)'
$ws.Range('L2').Value = 'maa://39402 (94.34), *maa://24633 (56.36), *maa://30515 (70.09), *maa://34787 (73.56), ***maa://20792 (11.93), ***maa://29083 (27.78)'
$ws.Range('T2').Value = 'maa://22742 (91.37), *maa://20791 (62.2)'
$ws.Range('AB2').Value = 'maa://21246 (91.42), maa://36684 (93.33), ***maa://22731 (6.25)'
$ws.Range('AF2').Value = 'maa://25251 (91.79), **maa://21730 (30.12), ***maa://39501 (10.87), **maa://36675 (50.0)'
$ws.Range('D3').Value = 'maa://40192 (96.77), maa://36987 (96.15), maa://39849 (88.89)'
$ws.Range('H3').Value = 'maa://21247 (98.39), *maa://22748 (60.0)'
$ws.Range('L3').Value = '*maa://22880 (65.33), maa://20276 (88.18), *maa://22749 (75.0)'
$ws.Range('P3').Value = 'maa://21249 (94.57), maa://26254 (97.14), **maa://22738 (50.0)'
$ws.Range('T3').Value = 'maa://24617 (90.62), maa://45854 (84.72), **maa://20790 (43.48), ***maa://37170 (16.67)'
$ws.Range('X3').Value = 'maa://27396 (84.1), maa://27484 (96.21), maa://27480 (83.78)'
$ws.Range('D4').Value = 'maa://24632 (93.81), **maa://24303 (38.46), maa://22499 (86.67), maa://22746 (100.0)'
$ws.Range('P4').Value = 'maa://49983 (94.74), maa://50121 (92.31)'
$ws.Range('T4').Value = 'maa://32509 (94.53), maa://27295 (88.51), maa://22754 (89.19), *maa://31008 (79.07), *maa://21746 (54.55)'
$ws.Range('X4').Value = '**maa://32495 (48.4), maa://43217 (92.98), ***maa://31785 (22.22), ***maa://36683 (29.79)'
$ws.Range('AF4').Value = '*maa://30062 (66.07), *maa://39394 (65.52), ***maa://26209 (13.04)'
$ws.Range('D5').Value = 'maa://21245 (84.78), maa://22744 (81.48)'
$ws.Range('P5').Value = 'maa://21919 (96.97), *maa://21281 (80.0)'
$ws.Range('AB5').Value = '*maa://29863 (63.41), ***maa://22752 (12.5), **maa://26013 (44.44)'
$ws.Range('D6').Value = 'maa://42407 (94.87)'
$ws.Range('X6').Value = 'maa://52754 (83.33)'
$ws.Range('AF6').Value = '*maa://33152 (62.71), ***maa://22770 (26.09)'
$ws.Range('P7').Value = 'maa://22750 (92.59)'
$ws.Range('T7').Value = 'maa://21291 (85.71)'
$ws.Range('X7').Value = 'maa://22399 (96.09), *maa://22758 (77.22)'
$ws.Range('AF7').Value = '*maa://26191 (67.82), maa://45272 (96.77), *maa://36671 (69.23), *maa://42530 (62.5)'
$ws.Range('A8').Value = '更新日期：2025.05.13 13:20:58'
$ws.Range('P8').Value = 'maa://32931 (84.33), *maa://21916 (62.69), maa://23252 (91.55), maa://37496 (97.3), **maa://22759 (45.45)'
$ws.Range('X8').Value = 'maa://21411 (96.14)'
$ws.Range('P9').Value = 'maa://22736 (83.04)'
$ws.Range('AB9').Value = 'maa://28711 (86.96), **maa://39938 (46.34), **maa://45044 (40.0), **maa://27377 (42.86), maa://40166 (97.14), ***maa://25174 (19.05)'
$ws.Range('AF9').Value = 'maa://26206 (88.41), *maa://22865 (51.85)'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = '6'
$ws.Range('D10').Value = '***maa://25695 (19.0), ***maa://39951 (12.5), ***maa://34206 (22.22), *maa://45271 (59.09), ***maa://39243 (25.0), maa://54000 (100.0)'
$ws.Range('T10').Value = 'maa://27395 (96.54), maa://22755 (88.89), **maa://22756 (40.91), ***maa://21737 (11.76)'
$ws.Range('X10').Value = 'maa://22301 (97.86), maa://45828 (92.86), maa://22726 (100.0)'
$ws.Range('L11').Value = 'maa://21287 (88.6)'
$ws.Range('T11').Value = 'maa://22747 (90.8), maa://22501 (98.11), maa://45521 (90.91)'
$ws.Range('X11').Value = 'maa://36713 (98.06)'
$ws.Range('AB11').Value = 'maa://29912 (97.7), maa://22516 (87.36), *maa://20794 (55.26)'
$ws.Range('AB12').Value = 'maa://23669 (95.64), maa://36677 (95.12), maa://39872 (93.1)'
$ws.Range('D13').Value = 'maa://24999 (92.58), maa://36673 (92.94), maa://25001 (86.3)'
$ws.Range('X13').Value = 'maa://34957 (82.35), **maa://22768 (50.0)'
$ws.Range('AF13').Value = '**maa://22737 (38.22), maa://39883 (87.07), *maa://39885 (51.28)'
$ws.Range('L14').Value = 'maa://39841 (94.0), maa://26245 (96.57), maa://21288 (96.3), maa://36682 (95.74)'
$ws.Range('D15').Value = '*maa://22743 (79.06), maa://22734 (84.43), *maa://30808 (65.22), *maa://36048 (65.66), maa://45058 (84.21)'
$ws.Range('H15').Value = 'maa://24304 (86.83), maa://21478 (90.0)'
$ws.Range('AF15').Value = 'maa://21364 (81.27), maa://36666 (81.25), *maa://22766 (68.29)'
$ws.Range('T16').Value = 'maa://22729 (94.44), *maa://28648 (72.73), *maa://36674 (79.31)'
$ws.Range('AF16').Value = '*maa://23911 (68.38), maa://27755 (93.88)'
$ws.Range('H17').Value = 'maa://22430 (89.47), maa://39599 (85.33)'
$ws.Range('D18').Value = 'maa://24570 (96.93)'
$ws.Range('H18').Value = 'maa://24421 (87.5)'
$ws.Range('AB18').Value = 'maa://24393 (98.11)'
$ws.Range('T19').Value = 'maa://24386 (98.61)'
$ws.Range('AB19').Value = '*maa://30709 (67.59), *maa://36668 (57.32)'
$ws.Range('D20').Value = 'maa://21432 (90.86), maa://25198 (94.44), *maa://20795 (50.76), maa://36680 (91.18)'
$ws.Range('L20').Value = 'maa://41331 (85.97)'
$ws.Range('X20').Value = 'maa://50085 (87.23), maa://49976 (85.9)'
$ws.Range('L21').Value = 'maa://31731 (96.49)'
$ws.Range('AB21').Value = 'maa://21443 (82.56), ***maa://23820 (30.0), **maa://52223 (46.67)'
$ws.Range('AF21').Value = 'maa://22524 (90.57), maa://22432 (82.3)'
$ws.Range('H22').Value = 'maa://25236 (96.12), **maa://21678 (48.94), **maa://22735 (42.86)'
$ws.Range('D23').Value = '***maa://28036 (26.92), **maa://41753 (48.0)'
$ws.Range('L23').Value = 'maa://39756 (95.44), maa://39875 (94.67)'
$ws.Range('P23').Value = 'maa://30587 (91.9), *maa://29748 (76.3), ***maa://29785 (17.39), *maa://37566 (78.26)'
$ws.Range('T23').Value = 'maa://31212 (91.67), maa://24387 (82.5)'
$ws.Range('X24').Value = 'maa://29988 (84.15), maa://23504 (93.56), **maa://22892 (41.29), *maa://25141 (77.37), *maa://36663 (77.78), ***maa://22815 (23.08), maa://52227 (100.0)'
$ws.Range('AF24').Value = 'maa://22523 (82.33), *maa://36672 (76.19), maa://29910 (93.75), **maa://21440 (35.71), maa://45831 (84.62)'
$ws.Range('D25').Value = 'maa://29753 (95.32)'
$ws.Range('H25').Value = '*maa://29063 (73.12), *maa://25311 (74.11), ***maa://22725 (4.76), *maa://45047 (72.73)'
$ws.Range('X25').Value = 'maa://29890 (82.76)'
$ws.Range('P26').Value = '*maa://30968 (65.0), maa://39870 (92.86)'
$ws.Range('H27').Value = '*maa://39601 (78.26), **maa://21283 (47.37), maa://34494 (97.14), **maa://36665 (50.0)'
$ws.Range('T27').Value = '*maa://30624 (77.14)'
$ws.Range('D28').Value = 'maa://24465 (90.61), maa://25725 (84.21)'
$ws.Range('T28').Value = '*maa://29765 (66.32), maa://23263 (95.41)'
$ws.Range('X28').Value = 'maa://39929 (91.82), maa://41749 (91.67), ***maa://39723 (13.89)'
$ws.Range('AF28').Value = 'maa://36660 (92.68), *maa://36701 (66.67)'
$ws.Range('H29').Value = '*maa://25175 (59.68)'
$ws.Range('AF29').Value = '*maa://24080 (69.17), maa://42865 (81.19), ***maa://34960 (8.33)'
$ws.Range('AB30').Value = 'maa://42979 (96.79), maa://45822 (100.0), maa://45045 (83.33)'
$ws.Range('L31').Value = 'maa://35926 (93.85), maa://36258 (86.93), *maa://43904 (71.43)'
$ws.Range('AB31').Value = '***maa://51420 (25.0)'
$ws.Range('H32').Value = 'maa://21895 (97.53), maa://36667 (98.02), **maa://20793 (38.0), maa://22760 (100.0)'
$ws.Range('T32').Value = 'maa://42859 (97.22), maa://41108 (88.0), maa://41238 (97.86), maa://45523 (100.0)'
$ws.Range('P34').Value = 'maa://48817 (96.23)'
$ws.Range('L35').Value = 'maa://41296 (97.13)'
$ws.Range('L37').Value = 'maa://45718 (98.3), *maa://47069 (78.95), maa://45789 (100.0)'
$ws.Range('H39').Value = 'maa://36670 (89.09), maa://25199 (85.22), maa://30434 (92.68), maa://45059 (80.65), ***maa://25036 (18.52), *maa://44165 (66.67)'
$ws.Range('T39').Value = 'maa://47079 (95.0), *maa://45788 (75.65), *maa://45790 (79.17)'
$ws.Range('P40').Value = 'maa://23278 (95.18), maa://21386 (95.79), maa://36664 (89.66), maa://45550 (87.5)'
$ws.Range('H44').Value = 'maa://29768 (98.23), maa://27728 (96.15)'
$ws.Range('H45').Value = 'maa://21229 (84.54), maa://30807 (94.67), *maa://22767 (55.0), maa://42459 (91.43), ***maa://20796 (13.79)'
$ws.Range('P45').Value = '*maa://36237 (75.0)'
$ws.Range('H46').Value = 'maa://35931 (92.29), maa://43901 (94.12)'
$ws.Range('H47').Value = 'maa://27410 (96.83), maa://29661 (97.55), maa://28038 (84.62)'
$ws.Range('H53').Value = 'maa://32534 (94.81), **maa://32434 (33.33)'
$ws.Range('H55').Value = 'maa://32532 (93.0)'
